$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -------------------------------------------------
$ws.Range("B1").Value = "Área Total [m²]"
$ws.Range("C1").Value = "Área Útil [m²]"

# Insert two new columns before the old "Localização" column (N) so the
# old N ("Localização") and O ("Link") shift right to P and Q.
$ws.Range("N1:O1").EntireColumn.Insert()
$ws.Range("N1").Value = "Matrícula"
$ws.Range("O1").Value = "Inscrição Imobiliária"

# --- Row 2 (Jose de Alencar property) -----------------------------------
$ws.Range("A2").Value = 209535.2
$ws.Range("B2").Value = 150
$ws.Range("C2").Value = 130
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 1
$ws.Range("H2").Value = 331170
$ws.Range("I2").Formula = "18/11/2021"
$ws.Range("M2").Formula = " Veras Negócios e Investimentos Ltda
"
$ws.Rows.Item(2).AutoFit()
$ws.Range("N2").Formula = "'62166"
$ws.Range("O2").Formula = "'5533031"

# --- Row 3 : replace with the "Coaçu" property --------------------------
$ws.Range("A3").Value = 720000
$ws.Range("B3").Value = 250
$ws.Range("C3").Value = 121
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = 720000
$ws.Range("G3").Value = 432000
$ws.Range("H3").Value = ""
$ws.Range("I3").Formula = "29/02/2024"
$ws.Range("K3").Formula = "'03/06/2024"
$ws.Range("L3").Formula = "'12/06/2024"
$ws.Range("M3").Formula = " Leilão Caixa "
$ws.Range("N3").Formula = "'73457"
$ws.Range("O3").Formula = "'5045363"
$ws.Range("P3").Formula = "RUA MARIA ALMEIDA ANTIGA RUA 03,N. 427 TERRENO 23A, COACU - CEP: 60871-742, FORTALEZA - CEARA"
$ws.Range("Q3").Formula = "https://www.leilaoimovel.com.br/imovel/ce/fortaleza/residencial-coacu-6-quartos-7-wc-wc-emp-4-salas-cozinha-imovel-caixa-economica-federal-cef-1620167-1188800008524-venda-direta-caixa"

# --- Remove the now-obsolete rows 4-6 -----------------------------------
$ws.Range("A4:A6").EntireRow.Delete()
